$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.611.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -7.35%  '
$ws.Range("D3").Value = "'1.697.54"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -6.00%  '
$ws.Range("D5").Value = "'219.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.35%  '
$ws.Range("D6").Value = "'0.5161"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -13.01%  '
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("E8").Value = '  -4.57%  '
$ws.Range("D9").Value = "'22.24"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.73%  '
$ws.Range("D10").Value = "'0.06276"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.03%  '
$ws.Range("D11").Value = "'0.07361"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.96%  '
$ws.Range("D12").Value = "'1.700.39"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -5.90%  '
$ws.Range("D13").Value = "'4.529"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.82%  '
$ws.Range("D14").Value = "'0.5860"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -6.01%  '
$ws.Range("D15").Value = "'1.929.64"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -5.91%  '
$ws.Range("D16").Value = "'0.000008412"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -9.08%  '
$ws.Range("D17").Value = "'65.68"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -13.19%  '
$ws.Range("D18").Value = "'26.668.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -7.00%  '
$ws.Range("D19").Value = "'5.039"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -8.01%  '
$ws.Range("E20").Value = '  +0.12%  '
$ws.Range("D21").Value = "'10.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.06%  '
$ws.Range("D22").Value = "'187.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -11.38%  '
$ws.Range("D23").Value = "'6.286"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -7.87%  '
$ws.Range("E24").Value = '  +0.23%  '
$ws.Range("D25").Value = "'144.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -5.85%  '
$ws.Range("D26").Value = "'7.609"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -3.42%  '
$ws.Range("D27").Value = "'0.1155"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -8.90%  '
$ws.Range("D28").Value = "'15.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.73%  '
$ws.Range("D29").Value = "'1.322"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -7.50%  '
$ws.Range("D30").Value = "'0.05700"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.46%  '
$ws.Range("D31").Value = "'1.343"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.97%  '
$ws.Range("D32").Value = "'3.523"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.98%  '
$ws.Range("D33").Value = "'3.518"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.15%  '
$ws.Range("D34").Value = "'1.660"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.21%  '
$ws.Range("D35").Value = "'1.027"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.28%  '
$ws.Range("D36").Value = "'0.6033"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -6.06%  '
$ws.Range("E37").Value = '  -5.01%  '
$ws.Range("D38").Value = "'2.689"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("D39").Value = "'1.102.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.74%  '
$ws.Range("D40").Value = "'0.01603"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.51%  '
$ws.Range("D41").Value = "'0.8626"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.17%  '
$ws.Range("D42").Value = "'5.874"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.84%  '
$ws.Range("E43").Value = '  -0.45%  '
$ws.Range("D44").Value = "'99.07"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.99%  '
$ws.Range("D45").Value = "'1.857.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.25%  '
$ws.Range("D46").Value = "'0.00000000109"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.17%  '
$ws.Range("D47").Value = "'56.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.22%  '
$ws.Range("D48").Value = "'8.175"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.65%  '
$ws.Range("D49").Value = "'1.003"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = "'0.05246"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.10%  '
$ws.Range("E51").Value = '  -3.40%  '
